$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.535.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.15%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.726.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.80%  "

# Row 4
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3734"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.32%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3399"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.94%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.183"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.75%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07447"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.56%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.416"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.74%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.53%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.054"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.722.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001074"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06665"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.74%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.194"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.511.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.99%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.453"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.430"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +21.61%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.392"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.43%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.916.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.107"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.992"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.26%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08624"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.64%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.694"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.362"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.07%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02343"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2161"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.90%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06215"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.427"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.74%  "

# Row 42
$ws.Range("E42").Value = "  -0.44%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6210"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.95%  "

# Row 44
$ws.Range("E44").Value = "  +6.47%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.891"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6019"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.46%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.043"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.95%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07177"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.08%  "
